$wb = $excel.ActiveWorkbook

# --- Shared string update: Tire_Type "710Rバフ100" -> "710R" ---
# Update the Tire_Type cell (K2:K6) on each Step3_DataPts_* sheet
foreach ($sheetName in @("Step3_DataPts_0.5", "Step3_DataPts_0.7", "Step3_DataPts_0.8", "Step3_DataPts_0.9")) {
    $s = $wb.Worksheets.Item($sheetName)
    foreach ($r in 2..6) {
        $s.Cells.Item($r, 11).Value = "710R"
    }
}

# --- Step1_Data updates ---
$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws1.Range("F2").Value = 0.0121840640905508
$ws1.Range("G2").Value = 0.08922939557640905
$ws1.Range("H2").Value = 0.2832256969257736
$ws1.Range("I2").Value = 0.05694755999134838
$ws1.Range("N2").Value = 0.06100340822144994
$ws1.Range("O2").Value = 0.009704892186060191
$ws1.Range("P2").Value = 0.1292570666351432
$ws1.Range("R2").Value = 0.05427572452530702
$ws1.Range("T2").Value = 0.166164902802697
$ws1.Range("U2").Value = 0.002108333739221937
$ws1.Range("V2").Value = 0.02065941416122502
$ws1.Range("X2").Value = 0.001108433611980476
$ws1.Range("Y2").Value = 0.008184638797427743
$ws1.Range("AA2").Value = 0.1012833523021285
$ws1.Range("AI2").Value = 0.00466311643327734
$ws1.Range("D3").Value = 0.1055350176767106
$ws1.Range("F3").Value = 0.3103978462646136
$ws1.Range("J3").Value = 0.00275752704865023
$ws1.Range("K3").Value = 0.01357331775101693
$ws1.Range("L3").Value = 0.0399150992472909
$ws1.Range("M3").Value = 0.07501213579482692
$ws1.Range("N3").Value = 0.06338371578751906
$ws1.Range("P3").Value = 0.07295185634979738
$ws1.Range("R3").Value = 0.1663957528291079
$ws1.Range("T3").Value = 0.01323194279838652
$ws1.Range("V3").Value = 0.0004537126269567694
$ws1.Range("W3").Value = 0.02166589749702399
$ws1.Range("Y3").Value = 0.1147261783280992
$ws1.Range("F4").Value = 0.2001272947229348
$ws1.Range("G4").Value = 0.1145502844178472
$ws1.Range("H4").Value = 0.1561787415873084
$ws1.Range("M4").Value = 0.05224032352961187
$ws1.Range("O4").Value = 0.1598342015039639
$ws1.Range("R4").Value = 0.0442379495556518
$ws1.Range("S4").Value = 0.08492486620232657
$ws1.Range("T4").Value = 0.06038485351835393
$ws1.Range("U4").Value = 0.01457506691445173
$ws1.Range("W4").Value = 0.009408114883704057
$ws1.Range("Y4").Value = 0.003039437697610227
$ws1.Range("Z4").Value = 0.06154909923823473
$ws1.Range("AA4").Value = 0.03806972138962724
$ws1.Range("AH4").Value = 0.0008800448383736487
$ws1.Range("F5").Value = 0.2210727623811943
$ws1.Range("G5").Value = 0.1022687849909022
$ws1.Range("H5").Value = 0.1792287689774957
$ws1.Range("J5").Value = 0.003998915462083324
$ws1.Range("M5").Value = 0.03581499511214403
$ws1.Range("O5").Value = 0.1572809883640796
$ws1.Range("R5").Value = 0.04371799225233787
$ws1.Range("S5").Value = 0.07988420072901989
$ws1.Range("T5").Value = 0.06859404301266761
$ws1.Range("U5").Value = 0.009669896348997901
$ws1.Range("W5").Value = 0.007404544147090041
$ws1.Range("Z5").Value = 0.06017408772060775
$ws1.Range("AA5").Value = 0.02683513440648752
$ws1.Range("AH5").Value = 0.004054886094891918
$ws1.Range("F6").Value = 0.08321948940234518
$ws1.Range("G6").Value = 0.01371716348530353
$ws1.Range("H6").Value = 0.3742390443889134
$ws1.Range("I6").Value = 0.003444437425260722
$ws1.Range("J6").Value = 0.03417921869716412
$ws1.Range("M6").Value = 0.01348852287713748
$ws1.Range("N6").Value = 0.01755558648863176
$ws1.Range("O6").Value = 0.08852257815998092
$ws1.Range("P6").Value = 0.0482915235685217
$ws1.Range("R6").Value = 0.07232279725165731
$ws1.Range("T6").Value = 0.1599966205821402
$ws1.Range("V6").Value = 0.01801966637521168
$ws1.Range("Y6").Value = 0.004177594177326839
$ws1.Range("Z6").Value = 0.01419644931112323
$ws1.Range("AA6").Value = 0.04956669128410442
$ws1.Range("AE6").Value = 0.002059178573404361
$ws1.Range("AI6").Value = 0.003003437951772948

# --- Step2_Sj updates (cumulative sums) ---
$ws2 = $wb.Worksheets.Item("Step2_Sj")
$ws2.Range("F2").Value = 0.0121840640905508
$ws2.Range("G2").Value = 0.10141345966695985
$ws2.Range("H2").Value = 0.38463915659273346
$ws2.Range("I2").Value = 0.44158671658408183
$ws2.Range("J2").Value = 0.44158671658408183
$ws2.Range("K2").Value = 0.44158671658408183
$ws2.Range("L2").Value = 0.44158671658408183
$ws2.Range("M2").Value = 0.44158671658408183
$ws2.Range("N2").Value = 0.5025901248055318
$ws2.Range("O2").Value = 0.512295016991592
$ws2.Range("P2").Value = 0.6415520836267352
$ws2.Range("Q2").Value = 0.6415520836267352
$ws2.Range("R2").Value = 0.6958278081520423
$ws2.Range("S2").Value = 0.6958278081520423
$ws2.Range("T2").Value = 0.8619927109547393
$ws2.Range("U2").Value = 0.8641010446939612
$ws2.Range("V2").Value = 0.8847604588551863
$ws2.Range("W2").Value = 0.8847604588551863
$ws2.Range("X2").Value = 0.8858688924671667
$ws2.Range("Y2").Value = 0.8940535312645944
$ws2.Range("Z2").Value = 0.8940535312645944
$ws2.Range("AA2").Value = 0.9953368835667229
$ws2.Range("AB2").Value = 0.9953368835667229
$ws2.Range("AC2").Value = 0.9953368835667229
$ws2.Range("AD2").Value = 0.9953368835667229
$ws2.Range("AE2").Value = 0.9953368835667229
$ws2.Range("AF2").Value = 0.9953368835667229
$ws2.Range("AG2").Value = 0.9953368835667229
$ws2.Range("AH2").Value = 0.9953368835667229
$ws2.Range("D3").Value = 0.1055350176767106
$ws2.Range("E3").Value = 0.1055350176767106
$ws2.Range("F3").Value = 0.41593286394132417
$ws2.Range("G3").Value = 0.41593286394132417
$ws2.Range("H3").Value = 0.41593286394132417
$ws2.Range("I3").Value = 0.41593286394132417
$ws2.Range("J3").Value = 0.4186903909899744
$ws2.Range("K3").Value = 0.43226370874099135
$ws2.Range("L3").Value = 0.47217880798828227
$ws2.Range("M3").Value = 0.5471909437831092
$ws2.Range("N3").Value = 0.6105746595706283
$ws2.Range("O3").Value = 0.6105746595706283
$ws2.Range("P3").Value = 0.6835265159204257
$ws2.Range("Q3").Value = 0.6835265159204257
$ws2.Range("R3").Value = 0.8499222687495336
$ws2.Range("S3").Value = 0.8499222687495336
$ws2.Range("T3").Value = 0.86315421154792
$ws2.Range("U3").Value = 0.86315421154792
$ws2.Range("V3").Value = 0.8636079241748769
$ws2.Range("W3").Value = 0.8852738216719008
$ws2.Range("X3").Value = 0.8852738216719008
$ws2.Range("F4").Value = 0.2001272947229348
$ws2.Range("G4").Value = 0.314677579140782
$ws2.Range("H4").Value = 0.47085632072809036
$ws2.Range("I4").Value = 0.47085632072809036
$ws2.Range("J4").Value = 0.47085632072809036
$ws2.Range("K4").Value = 0.47085632072809036
$ws2.Range("L4").Value = 0.47085632072809036
$ws2.Range("M4").Value = 0.5230966442577022
$ws2.Range("N4").Value = 0.5230966442577022
$ws2.Range("O4").Value = 0.6829308457616661
$ws2.Range("P4").Value = 0.6829308457616661
$ws2.Range("Q4").Value = 0.6829308457616661
$ws2.Range("R4").Value = 0.7271687953173178
$ws2.Range("S4").Value = 0.8120936615196443
$ws2.Range("T4").Value = 0.8724785150379983
$ws2.Range("U4").Value = 0.8870535819524501
$ws2.Range("V4").Value = 0.8870535819524501
$ws2.Range("W4").Value = 0.8964616968361542
$ws2.Range("X4").Value = 0.8964616968361542
$ws2.Range("Y4").Value = 0.8995011345337645
$ws2.Range("Z4").Value = 0.9610502337719992
$ws2.Range("AA4").Value = 0.9991199551616264
$ws2.Range("AB4").Value = 0.9991199551616264
$ws2.Range("AC4").Value = 0.9991199551616264
$ws2.Range("AD4").Value = 0.9991199551616264
$ws2.Range("AE4").Value = 0.9991199551616264
$ws2.Range("AF4").Value = 0.9991199551616264
$ws2.Range("AG4").Value = 0.9991199551616264
$ws2.Range("F5").Value = 0.2210727623811943
$ws2.Range("G5").Value = 0.3233415473720965
$ws2.Range("H5").Value = 0.5025703163495923
$ws2.Range("I5").Value = 0.5025703163495923
$ws2.Range("J5").Value = 0.5065692318116756
$ws2.Range("K5").Value = 0.5065692318116756
$ws2.Range("L5").Value = 0.5065692318116756
$ws2.Range("M5").Value = 0.5423842269238196
$ws2.Range("N5").Value = 0.5423842269238196
$ws2.Range("O5").Value = 0.6996652152878993
$ws2.Range("P5").Value = 0.6996652152878993
$ws2.Range("Q5").Value = 0.6996652152878993
$ws2.Range("R5").Value = 0.7433832075402371
$ws2.Range("S5").Value = 0.823267408269257
$ws2.Range("T5").Value = 0.8918614512819246
$ws2.Range("U5").Value = 0.9015313476309226
$ws2.Range("V5").Value = 0.9015313476309226
$ws2.Range("W5").Value = 0.9089358917780126
$ws2.Range("X5").Value = 0.9089358917780126
$ws2.Range("Y5").Value = 0.9089358917780126
$ws2.Range("Z5").Value = 0.9691099794986203
$ws2.Range("AA5").Value = 0.9959451139051079
$ws2.Range("AB5").Value = 0.9959451139051079
$ws2.Range("AC5").Value = 0.9959451139051079
$ws2.Range("AD5").Value = 0.9959451139051079
$ws2.Range("AE5").Value = 0.9959451139051079
$ws2.Range("AF5").Value = 0.9959451139051079
$ws2.Range("AG5").Value = 0.9959451139051079
$ws2.Range("F6").Value = 0.08321948940234518
$ws2.Range("G6").Value = 0.0969366528876487
$ws2.Range("H6").Value = 0.4711756972765621
$ws2.Range("I6").Value = 0.47462013470182285
$ws2.Range("J6").Value = 0.508799353398987
$ws2.Range("K6").Value = 0.508799353398987
$ws2.Range("L6").Value = 0.508799353398987
$ws2.Range("M6").Value = 0.5222878762761245
$ws2.Range("N6").Value = 0.5398434627647563
$ws2.Range("O6").Value = 0.6283660409247371
$ws2.Range("P6").Value = 0.6766575644932589
$ws2.Range("Q6").Value = 0.6766575644932589
$ws2.Range("R6").Value = 0.7489803617449162
$ws2.Range("S6").Value = 0.7489803617449162
$ws2.Range("T6").Value = 0.9089769823270564
$ws2.Range("U6").Value = 0.9089769823270564
$ws2.Range("V6").Value = 0.9269966487022681
$ws2.Range("W6").Value = 0.9269966487022681
$ws2.Range("X6").Value = 0.9269966487022681
$ws2.Range("Y6").Value = 0.9311742428795949
$ws2.Range("Z6").Value = 0.9453706921907181
$ws2.Range("AA6").Value = 0.9949373834748225
$ws2.Range("AB6").Value = 0.9949373834748225
$ws2.Range("AC6").Value = 0.9949373834748225
$ws2.Range("AD6").Value = 0.9949373834748225
$ws2.Range("AE6").Value = 0.9969965620482268
$ws2.Range("AF6").Value = 0.9969965620482268
$ws2.Range("AG6").Value = 0.9969965620482268
$ws2.Range("AH6").Value = 0.9969965620482268

# --- Step3_DataPts_0.5 updates ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F2").Value = 0.5025901248055318
$ws.Range("D3").Value = 12
$ws.Range("F3").Value = 0.5471909437831092
$ws.Range("G3").Value = 11
$ws.Range("D4").Value = 12
$ws.Range("F4").Value = 0.5230966442577022
$ws.Range("G4").Value = 9
$ws.Range("F5").Value = 0.5025703163495923
$ws.Range("D6").Value = 9
$ws.Range("F6").Value = 0.508799353398987
$ws.Range("G6").Value = 6

# --- Step3_DataPts_0.7 updates ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 19
$ws.Range("F2").Value = 0.8619927109547393
$ws.Range("G2").Value = 15
$ws.Range("F3").Value = 0.8499222687495336
$ws.Range("D4").Value = 17
$ws.Range("F4").Value = 0.7271687953173178
$ws.Range("G4").Value = 14
$ws.Range("D5").Value = 17
$ws.Range("F5").Value = 0.7433832075402371
$ws.Range("G5").Value = 14
$ws.Range("D6").Value = 17
$ws.Range("F6").Value = 0.7489803617449162
$ws.Range("G6").Value = 14

# --- Step3_DataPts_0.8 updates ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F2").Value = 0.8619927109547393
$ws.Range("F3").Value = 0.8499222687495336
$ws.Range("D4").Value = 18
$ws.Range("F4").Value = 0.8120936615196443
$ws.Range("G4").Value = 15
$ws.Range("D5").Value = 18
$ws.Range("F5").Value = 0.823267408269257
$ws.Range("G5").Value = 15
$ws.Range("F6").Value = 0.9089769823270564

# --- Step3_DataPts_0.9 updates ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F2").Value = 0.9953368835667229
$ws.Range("D4").Value = 25
$ws.Range("F4").Value = 0.9610502337719992
$ws.Range("G4").Value = 22
$ws.Range("D5").Value = 20
$ws.Range("F5").Value = 0.9015313476309226
$ws.Range("G5").Value = 17
$ws.Range("F6").Value = 0.9089769823270564
